$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '58.637.06'
$ws.Cells.Item(2, 5).Value = '  -1.40%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.628.11'
$ws.Cells.Item(3, 5).Value = '  +0.87%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '534.44'
$ws.Cells.Item(5, 5).Value = '  -0.71%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '142.73'
$ws.Cells.Item(6, 5).Value = '  +0.78%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.998'
$ws.Cells.Item(7, 5).Value = '  -0.08%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.568'
$ws.Cells.Item(8, 5).Value = '  +0.32%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '2.635.54'
$ws.Cells.Item(9, 5).Value = '  +1.00%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '7.02'
$ws.Cells.Item(10, 5).Value = '  +8.84%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.101'
$ws.Cells.Item(11, 5).Value = '  -1.67%  '
$ws.Cells.Item(12, 5).Value = '  +0.08%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.135'
$ws.Cells.Item(13, 5).Value = '  +0.67%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '3.090.67'
$ws.Cells.Item(14, 5).Value = '  +0.92%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '58.587.90'
$ws.Cells.Item(15, 5).Value = '  -1.34%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '20.91'
$ws.Cells.Item(16, 5).Value = '  +1.27%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.625.75'
$ws.Cells.Item(17, 5).Value = '  +0.49%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.0000133'
$ws.Cells.Item(18, 5).Value = '  -0.58%  '
$ws.Cells.Item(19, 2).Value = 'Polkadot'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '4.40'
$ws.Cells.Item(19, 5).Value = '  +0.86%  '
$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '334.51'
$ws.Cells.Item(20, 5).Value = '  -2.19%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '10.17'
$ws.Cells.Item(21, 5).Value = '  +0.66%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.21'
$ws.Cells.Item(23, 5).Value = '  -0.10%  '
$ws.Cells.Item(24, 5).Value = '  -1.48%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.417'
$ws.Cells.Item(25, 5).Value = '  +2.01%  '
$ws.Cells.Item(26, 5).Value = '  -0.54%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.998'
$ws.Cells.Item(27, 5).Value = '  -0.06%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.13'
$ws.Cells.Item(28, 5).Value = '  -1.10%  '
$ws.Cells.Item(29, 2).Value = 'PEPE'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.0₃0737'
$ws.Cells.Item(29, 5).Value = '  -1.33%  '
$ws.Cells.Item(30, 2).Value = 'USDe'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.998'
$ws.Cells.Item(30, 5).Value = '  -0.11%  '
$ws.Cells.Item(31, 5).Value = '  -2.02%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '5.89'
$ws.Cells.Item(32, 5).Value = '  +1.04%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '18.77'
$ws.Cells.Item(33, 5).Value = '  -0.39%  '
$ws.Cells.Item(34, 2).Value = 'Monero'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '150.59'
$ws.Cells.Item(34, 5).Value = '  +0.45%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '3.91'
$ws.Cells.Item(35, 5).Value = '  -1.70%  '
$ws.Cells.Item(36, 2).Value = 'OKB'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '37.09'
$ws.Cells.Item(36, 5).Value = '  -0.12%  '
$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.11'
$ws.Cells.Item(37, 5).Value = '  -0.76%  '
$ws.Cells.Item(38, 2).Value = 'SuiNetwork'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.824'
$ws.Cells.Item(38, 5).Value = '  -2.56%  '
$ws.Cells.Item(39, 2).Value = 'Stacks'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '1.43'
$ws.Cells.Item(39, 5).Value = '  -2.70%  '
$ws.Cells.Item(40, 2).Value = 'Fetch.AI'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.814'
$ws.Cells.Item(40, 5).Value = '  -1.61%  '
$ws.Cells.Item(41, 2).Value = 'Filecoin'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.59'
$ws.Cells.Item(41, 5).Value = '  +1.38%  '
$ws.Cells.Item(42, 2).Value = 'Bittensor'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '281.41'
$ws.Cells.Item(42, 5).Value = '  +3.11%  '
$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.999'
$ws.Cells.Item(43, 5).Value = '  -0.07%  '
$ws.Cells.Item(44, 2).Value = 'Mantle'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.600'
$ws.Cells.Item(44, 5).Value = '  +0.66%  '
$ws.Cells.Item(45, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '10.70'
$ws.Cells.Item(45, 5).Value = '  -0.27%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '19.17'
$ws.Cells.Item(46, 5).Value = '  +3.45%  '
$ws.Cells.Item(47, 2).Value = 'Hedera'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.0532'
$ws.Cells.Item(47, 5).Value = '  +1.76%  '
$ws.Cells.Item(48, 2).Value = 'Stellar'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.0936'
$ws.Cells.Item(48, 5).Value = '  -2.03%  '
$ws.Cells.Item(49, 2).Value = 'VeChain'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.0225'
$ws.Cells.Item(49, 5).Value = '  +0.86%  '
$ws.Cells.Item(50, 2).Value = 'Maker'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.941.34'
$ws.Cells.Item(50, 5).Value = '  -0.41%  '
$ws.Cells.Item(51, 2).Value = 'RenderToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '4.46'
$ws.Cells.Item(51, 5).Value = '  -1.25%  '
